# PopulationCalculator.xlsx - "Updating demographics & religion"
#
# The workbook tracks population-by-race percentages for each settlement.
# Each settlement row (e.g. Stilben, Rural 50% of Stilben, Rural 200%
# Bronbog) carries one input percentage per race column; a hidden "sum
# check" column verifies the row's percentages add up to 100%. This edit
# rebalances a handful of race percentages on three settlement rows while
# keeping every row summing to exactly 1 (100%). All the downstream
# per-race head-count / percentage-of-total formulas recalculate
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30 : Stilben ------------------------------------------------
$ws.Range("H30").Value  = 0.0175   # Aasimar      19.0% -> 17.5%
$ws.Range("N30").Value  = 0.019    # Dragonborn    2.0% -> 1.9%
$ws.Range("P30").Value  = 0.001    # Drow          0.0% -> 0.1%
$ws.Range("BL30").Value = 0.0025   # Yuan-Ti       0.1% -> 0.25%

# --- Row 32 : Rural (50% of Stilben) ----------------------------------
$ws.Range("J32").Value  = 0.0225   # Bugbear       3.0% -> 2.25%
$ws.Range("P32").Value  = 0        # Drow          0.1% -> 0.0%
$ws.Range("V32").Value  = 0.025    # Elf           4.0% -> 2.5%
$ws.Range("AB32").Value = 0.075    # Gnoll         3.0% -> 7.5%
$ws.Range("AD32").Value = 0.015    # Gnome         2.0% -> 1.5%
$ws.Range("AF32").Value = 0.035    # Goblin        2.0% -> 3.5%
$ws.Range("AJ32").Value = 0.025    # Half Elf      4.8% -> 2.5%
$ws.Range("AL32").Value = 0.04     # Half Orc      5.5% -> 4.0%
$ws.Range("AN32").Value = 0.035    # Halfling      4.5% -> 3.5%
$ws.Range("AR32").Value = 0.37     # Human        41.0% -> 37.0%
$ws.Range("AV32").Value = 0.02     # Kobold        3.0% -> 2.0%
$ws.Range("AX32").Value = 0.105    # Lizardfolk    6.0% -> 10.5%
$ws.Range("BB32").Value = 0.11     # Orc           5.0% -> 11.0%
$ws.Range("BF32").Value = 0.01     # Tieflings     5.0% -> 1.0%
$ws.Range("BL32").Value = 0.0025   # Yuan-Ti       0.1% -> 0.25%

# --- Row 45 : Rural (200% Bronbog) ------------------------------------
$ws.Range("V45").Value  = 0.015    # Elf           2.0% -> 1.5%
$ws.Range("AB45").Value = 0.05     # Gnoll         4.0% -> 5.0%
$ws.Range("AD45").Value = 0.005    # Gnome         1.0% -> 0.5%
$ws.Range("AF45").Value = 0.11     # Goblin       10.0% -> 11.0%
$ws.Range("AH45").Value = 0.02     # Goliath       3.0% -> 2.0%
$ws.Range("AJ45").Value = 0.04     # Half Elf      9.0% -> 4.0%
$ws.Range("AL45").Value = 0.005    # Half Orc      0.0% -> 0.5%
$ws.Range("AR45").Value = 0.24     # Human        32.0% -> 24.0%
$ws.Range("AV45").Value = 0.13     # Kobold       11.0% -> 13.0%
$ws.Range("AX45").Value = 0.11     # Lizardfolk    8.0% -> 11.0%
$ws.Range("BB45").Value = 0.14     # Orc           9.0% -> 14.0%
$ws.Range("BF45").Value = 0.005    # Tieflings     0.0% -> 0.5%
$ws.Range("BH45").Value = 0.04     # Tortle        2.0% -> 4.0%

# --- View state: scroll/zoom/selection as left by the editor ---------
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("AD126").Select()
